$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metrics for row 20 (2025Q2): total_customers, returning_customers,
# new_customers and recurrence_rate were refreshed with the new data pull.
$ws.Range("C20").Value = 321
$ws.Range("D20").Value = 257
$ws.Range("E20").Value = 64
$ws.Range("F20").Value = 79.56656346749226
